$d = $word.ActiveDocument

# Locate the paragraph that still needs to be crossed out:
# "Connect Interface able to connect to at least one online shop"
$range = $d.Content
$found = $range.Find.Execute("Connect Interface able to connect to at least one online shop", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Expand the range to include the paragraph mark so the paragraph-mark
    # run properties (pPr/rPr) also get the strikethrough, matching Word's
    # behavior when you select the whole line (e.g. triple-click) and
    # press the Strikethrough button.
    $para = $range.Paragraphs(1)
    $paraRange = $para.Range
    $paraRange.Font.StrikeThrough = $true
}
